# ---------------------------------------------------------------------------
# Split the single "Generat token:<token>" paragraph into two paragraphs:
#   1) "Generat token:" + the original token (bookmark _GoBack removed here)
#   2) the new token, with the _GoBack bookmark now anchored at its end
# Both paragraphs also pick up an explicit east-Asia font hint on the
# paragraph mark itself (w:pPr/w:rPr/w:rFonts[@w:hint='eastAsia']).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

$targetParagraph = $d.Paragraphs(1)

$newXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"
                   xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
        <w:body>
          <w:p w14:paraId="464B2D33">
            <w:pPr>
              <w:ind w:firstLine="210" w:firstLineChars="100"/>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
                <w:lang w:val="en-US" w:eastAsia="zh-CN"/>
              </w:rPr>
              <w:t>Generat token:</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
              </w:rPr>
              <w:t>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</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:ind w:firstLine="210" w:firstLineChars="100"/>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:hint="eastAsia"/>
              </w:rPr>
              <w:t>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</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$targetParagraph.Range.InsertXML($newXml)

# ---------------------------------------------------------------------------
# Mark the "Default Paragraph Font" character style as a Quick Style
# (adds <w:qFormat/> to its style definition in styles.xml).
# ---------------------------------------------------------------------------
$defaultParaFont = $d.Styles("Default Paragraph Font")
$defaultParaFont.QuickStyle = $true
